$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 35: bump Manhours and extend the Comment with the follow-up note
$ws.Range("B35").Value = 0.5
$ws.Range("F35").Value = "Once you’ve done one…then you spot bugs"

# New rows 36-43: additional "Tag budget" work items
$newRows = 36..43
foreach ($r in $newRows) {
    $ws.Cells.Item($r, 1).Value = $ws.Cells.Item(35, 1).Value2
    $ws.Cells.Item($r, 1).NumberFormat = $ws.Cells.Item(35, 1).NumberFormat
    $ws.Cells.Item($r, 3).Value = "Implementation"
    $ws.Cells.Item($r, 4).Value = 2
}

$ws.Range("E36").Value = "Add budget to tags table, controller, specs and views"

# Move the active selection to where the new content was entered
$ws.Range("E36").Select()
